# Insert a new price-record row at row 146 (shifts existing rows 146-250
# down to 147-251) and populate it with the new record's data.
# This mirrors the OOXML diff: dimension A1:R250 -> A1:R251, and every
# row from the old 146 through 250 reappears one row lower, with the
# brand-new row occupying the vacated row 146.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(146).Insert()

$ws.Range("A146").Value = 10
$ws.Range("B146").Value = "Vega Modelo de Temuco"
$ws.Range("C146").Value = "La Araucanía"
$ws.Range("D146").Value = 44574
$ws.Range("E146").Value = 9
$ws.Range("F146").Value = 100114013
$ws.Range("G146").Value = "Zanahoria"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 280
$ws.Range("K146").Value = 8000
$ws.Range("L146").Value = 8000
$ws.Range("M146").Value = 8000
$ws.Range("N146").Value = "$/saco 20 kilos"
$ws.Range("O146").Value = "Región del Maule"
$ws.Range("P146").Value = 400
$ws.Range("Q146").Value = 20
$ws.Range("R146").Value = "Hortaliza"
